# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 366
$ws1.Range("F3").Value = 778
$ws1.Range("F4").Value = 272
$ws1.Range("F5").Value = 834
$ws1.Range("F6").Value = 2060
$ws1.Range("F7").Value = 185

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 366
$ws4.Range("F3").Value = 778
$ws4.Range("F4").Value = 272
$ws4.Range("F7").Value = 834
$ws4.Range("F8").Value = 2060
$ws4.Range("F10").Value = 185
